$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")

$ws.Range("T1").Value = "Year"

$ws.Range("T2").Value = 2010
$ws.Range("T3").Value = 2010
$ws.Range("T4").Value = 2020
$ws.Range("T5").Value = 2024
$ws.Range("T6").Value = 2025
$ws.Range("T7").Value = 2019
$ws.Range("T8").Value = 2018
$ws.Range("T9").Value = 2010
$ws.Range("T10").Value = 2019
$ws.Range("T11").Value = 2010
$ws.Range("T12").Value = 2010
$ws.Range("T13").Value = 2020
$ws.Range("T14").Value = 2024
$ws.Range("T15").Value = 2025
$ws.Range("T16").Value = 2019
$ws.Range("T17").Value = 2018
$ws.Range("T18").Value = 2010
$ws.Range("T19").Value = 2019
$ws.Range("T20").Value = 2025
$ws.Range("T21").Value = 2019
$ws.Range("T22").Value = 2018
$ws.Range("T23").Value = 2010
$ws.Range("T24").Value = 2019
$ws.Range("T25").Value = 2015
$ws.Range("T26").Value = 2010
$ws.Range("T27").Value = 2010
$ws.Range("T28").Value = 2020
$ws.Range("T29").Value = 2024
$ws.Range("T30").Value = 2025
$ws.Range("T31").Value = 2019
$ws.Range("T32").Value = 2018
$ws.Range("T33").Value = 2010
$ws.Range("T34").Value = 2019
$ws.Range("T35").Value = 2016
$ws.Range("T36").Value = 2018
$ws.Range("T37").Value = 2014
$ws.Range("T38").Value = 2016
$ws.Range("T39").Value = 2017
$ws.Range("T40").Value = 2019
$ws.Range("T41").Value = 2020
$ws.Range("T42").Value = 2024
$ws.Range("T43").Value = 2021
$ws.Range("T44").Value = 2023
$ws.Range("T45").Value = 2022
$ws.Range("T46").Value = 2020
$ws.Range("T47").Value = 2021
$ws.Range("T48").Value = 2025
$ws.Range("T49").Value = 2024
$ws.Range("T50").Value = 2024
$ws.Range("T51").Value = 2020
